# Auto-generated edit script applying the Rafflesia_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 324.75
$ws.Range("I29").Value = 324.75
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 974.25
$ws.Range("L29").Value = 0
$ws.Range("M29").Value = -693.25
$ws.Range("N29").Value = ""

$ws.Range("H31").Value = 450.25
$ws.Range("I31").Value = 450.25
$ws.Range("K31").Value = 1350.75
$ws.Range("M31").Value = -1120.75

$ws.Range("H39").Value = 23.75
$ws.Range("I39").Value = 23.75
$ws.Range("K39").Value = 71.25
$ws.Range("M39").Value = 224.75

$ws.Range("H70").Value = 1124.5
$ws.Range("I70").Value = 899
$ws.Range("J70").Value = 1350
$ws.Range("K70").Value = 2697
$ws.Range("L70").Value = 4050
$ws.Range("M70").Value = -2427
$ws.Range("N70").Value = -4590

$ws.Range("H73").Value = 1124.5
$ws.Range("I73").Value = 899
$ws.Range("J73").Value = 1350
$ws.Range("K73").Value = 2697
$ws.Range("L73").Value = 4050
$ws.Range("M73").Value = -1761
$ws.Range("N73").Value = -5922

$ws.Range("H132").Value = 5553.2354
$ws.Range("I132").Value = 5622.643
$ws.Range("J132").Value = 5229.3335
$ws.Range("K132").Value = 16867.929
$ws.Range("L132").Value = 15688.0005
$ws.Range("M132").Value = -14337.929
$ws.Range("N132").Value = -20748.0005

$ws.Range("H137").Value = 1593.5
$ws.Range("I137").Value = 1593.5
$ws.Range("K137").Value = 4780.5
$ws.Range("M137").Value = -2230.5

$ws.Range("H138").Value = 197
$ws.Range("I138").Value = 197
$ws.Range("K138").Value = 591
$ws.Range("M138").Value = 4549

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1166
$ws.Range("I2").Value = 1166
$ws.Range("K2").Value = 1166
$ws.Range("M2").Value = -1053

$ws.Range("H32").Value = 6444.4546
$ws.Range("I32").Value = 6444.4546
$ws.Range("K32").Value = 6444.4546
$ws.Range("M32").Value = -6157.4546

$ws.Range("H45").Value = 2305.4167
$ws.Range("I45").Value = 2471.375
$ws.Range("J45").Value = 1973.5
$ws.Range("K45").Value = 2471.375
$ws.Range("L45").Value = 1973.5
$ws.Range("M45").Value = -2094.375
$ws.Range("N45").Value = -2727.5

$ws.Range("H116").Value = 1166
$ws.Range("I116").Value = 1166
$ws.Range("K116").Value = 1166
$ws.Range("M116").Value = 1128

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1166
$ws.Range("I3").Value = 1166
$ws.Range("K3").Value = 1166
$ws.Range("M3").Value = -1052

$ws.Range("H5").Value = 826.2857
$ws.Range("I5").Value = 401.33334
$ws.Range("J5").Value = 1145
$ws.Range("K5").Value = 401.33334
$ws.Range("L5").Value = 1145
$ws.Range("M5").Value = -288.33334
$ws.Range("N5").Value = -1371

$ws.Range("H7").Value = 87500170
$ws.Range("I7").Value = 87500170
$ws.Range("K7").Value = 87500170
$ws.Range("M7").Value = -87500057

$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").Value = ""

$ws.Range("H11").Value = 4566
$ws.Range("I11").Value = 300
$ws.Range("J11").Value = 7410
$ws.Range("K11").Value = 300
$ws.Range("L11").Value = 7410
$ws.Range("M11").Value = -160
$ws.Range("N11").Value = -7690

$ws.Range("H22").Value = 1500
$ws.Range("I22").Value = 1500
$ws.Range("K22").Value = 1500
$ws.Range("M22").Value = -1327

$ws.Range("H86").Value = 4069.625
$ws.Range("I86").Value = 4069.625
$ws.Range("K86").Value = 4069.625
$ws.Range("M86").Value = -2946.625

$ws.Range("H89").Value = 4069.625
$ws.Range("I89").Value = 4069.625
$ws.Range("K89").Value = 20348.125
$ws.Range("M89").Value = -14732.125

$ws.Range("H99").Value = 1069.875
$ws.Range("I99").Value = 1008.4286
$ws.Range("K99").Value = 1008.4286
$ws.Range("M99").Value = 489.5714

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 2482.25
$ws.Range("I2").Value = 3234.6667
$ws.Range("J2").Value = 225
$ws.Range("K2").Value = 3234.6667
$ws.Range("L2").Value = 225
$ws.Range("M2").Value = -3121.6667
$ws.Range("N2").Value = -451

$ws.Range("H31").Value = 2639.4546
$ws.Range("I31").Value = 1824.875
$ws.Range("J31").Value = 4811.6665
$ws.Range("K31").Value = 1824.875
$ws.Range("L31").Value = 4811.6665
$ws.Range("M31").Value = -1529.875
$ws.Range("N31").Value = -5401.6665

$ws.Range("H34").Value = 2639.4546
$ws.Range("I34").Value = 1824.875
$ws.Range("J34").Value = 4811.6665
$ws.Range("K34").Value = 1824.875
$ws.Range("L34").Value = 4811.6665
$ws.Range("M34").Value = -1622.875
$ws.Range("N34").Value = -5215.6665

$ws.Range("H122").Value = 980.8
$ws.Range("I122").Value = 985
$ws.Range("J122").Value = 964
$ws.Range("K122").Value = 2955
$ws.Range("L122").Value = 2892
$ws.Range("M122").Value = -505
$ws.Range("N122").Value = -7792

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 908.25
$ws.Range("J34").Value = 1214.1428
$ws.Range("L34").Value = 3642.4284
$ws.Range("N34").Value = -3810.4284

$ws.Range("H114").Value = 883.6667
$ws.Range("J114").Value = 825.5
$ws.Range("L114").Value = 2476.5
$ws.Range("N114").Value = -8984.5

$ws.Range("H124").Value = 5000
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = ""

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2148.8
$ws.Range("I40").Value = 2148.8
$ws.Range("K40").Value = 2148.8
$ws.Range("M40").Value = -2012.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 10000
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 10000
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 10000
$ws.Range("M63").Value = ""
$ws.Range("N63").Value = -11248

$ws.Range("H66").Value = 10000
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 10000
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 30000
$ws.Range("M66").Value = ""
$ws.Range("N66").Value = -36240

$ws.Range("H81").Value = 833.3333
$ws.Range("J81").Value = 1500
$ws.Range("L81").Value = 3000
$ws.Range("N81").Value = -5122

$ws.Range("H84").Value = 833.3333
$ws.Range("J84").Value = 1500
$ws.Range("L84").Value = 15000
$ws.Range("N84").Value = -25608

$ws.Range("H86").Value = 120000
$ws.Range("J86").Value = 120000
$ws.Range("L86").Value = 120000
$ws.Range("N86").Value = -122246

$ws.Range("H89").Value = 120000
$ws.Range("J89").Value = 120000
$ws.Range("L89").Value = 600000
$ws.Range("N89").Value = -611232

$ws.Range("H113").Value = 1642.7142
$ws.Range("I113").Value = 1050
$ws.Range("J113").Value = 1879.8
$ws.Range("K113").Value = 3150
$ws.Range("L113").Value = 5639.4
$ws.Range("M113").Value = -980
$ws.Range("N113").Value = -9979.4
